# Updated capital structure database
# Applies the capital-structure refresh to the Kazakhstan Metals & Mining
# rows (2 and 3): adds the new historical_growth_revenue_last_5_years
# figure (column D), recomputes the margin / cash-return / return / debt
# metrics, and drops the now-unused buybacks_cash_returned value
# (column T, which collapses into the buybacks column S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2, 3) {
    # New column: historical_growth_revenue_last_5_years
    $ws.Cells.Item($row, 4).Value = -0.0946          # D

    # Margin metrics
    $ws.Cells.Item($row, 7).Value  = -8.983333333333334   # G  ebitdard_margin
    $ws.Cells.Item($row, 8).Value  = -9.944444444444445   # H  ebitda_margin
    $ws.Cells.Item($row, 9).Value  = -9.722222222222223   # I  operating_margin
    $ws.Cells.Item($row, 10).Value = -9.722222222222223   # J  after_tax_operating_margin
    $ws.Cells.Item($row, 11).Value = -4.83                # K  trailing_net_income
    $ws.Cells.Item($row, 12).Value = -26.83333333333334   # L  net_margin

    # Cash returned
    $ws.Cells.Item($row, 13).Value = 0                    # M  cash_returned
    $ws.Cells.Item($row, 14).Value = 0                    # N  cash_returned_market_cap
    $ws.Cells.Item($row, 15).Value = 0                    # O  cash_returned_net_income

    # Buybacks: buybacks_cash_returned (T) is retired, buybacks (S) recomputed
    $ws.Cells.Item($row, 19).Value = 0                    # S  buybacks
    $ws.Cells.Item($row, 20).ClearContents()              # T  buybacks_cash_returned (removed)

    # Cash position
    $ws.Cells.Item($row, 21).Value = 0.018                     # U  cash
    $ws.Cells.Item($row, 22).Value = 0.0009137055837563451     # V  cash_market_cap

    # Returns on equity / capital
    $ws.Cells.Item($row, 23).Value = 0.6580381471389646    # W  roe
    $ws.Cells.Item($row, 24).Value = 0.111545147628466     # X  cost_equity
    $ws.Cells.Item($row, 25).Value = 0.5464929995104986    # Y  roe_cost_equity
    $ws.Cells.Item($row, 26).Value = 0.01539514197742045   # Z  sales_invested_capital
    $ws.Cells.Item($row, 27).Value = -0.1496749914471433   # AA roic
    $ws.Cells.Item($row, 28).Value = 0.09372376010761474   # AB cost_capital
    $ws.Cells.Item($row, 29).Value = -0.2433987515547581   # AC roic_cost_capital

    # Debt structure
    $ws.Cells.Item($row, 30).Value = 19.7                  # AD debt_total
    $ws.Cells.Item($row, 32).Value = 19.7                  # AF debt_total_inc_leases
    $ws.Cells.Item($row, 33).Value = 19.682                # AG net_debt
    $ws.Cells.Item($row, 34).Value = 0.5                   # AH debt_market_capital
    $ws.Cells.Item($row, 35).Value = 1.669491525423729     # AI debt_book_capital
    $ws.Cells.Item($row, 36).Value = 0.4997714691991265    # AJ net_debt_market_capital
    $ws.Cells.Item($row, 37).Value = 1.670514343914446     # AK net_debt_book_capital
    $ws.Cells.Item($row, 38).Value = 2.67                  # AL interest_expenses
    $ws.Cells.Item($row, 39).Value = 2.67                  # AM net_interest_expenses

    # Leverage / coverage ratios
    $ws.Cells.Item($row, 40).Value = -11.32183908045977    # AN debt_ebitda
    $ws.Cells.Item($row, 41).Value = -0.6554307116104869   # AO ebit_interest_expenses
    $ws.Cells.Item($row, 42).Value = -11.31149425287356    # AP net_debt_ebitda
    $ws.Cells.Item($row, 43).Value = -0.6554307116104869   # AQ ebit_net_interest_expenses
}

Write-Output "Kazakhstan metals & mining capital structure rows updated"
